$d = $word.ActiveDocument

# Helper: split the run(s) covering a given sub-string (identified via Find,
# scoped to a bounding range) into its own run by toggling a character
# formatting property on and back off. Word (and this COM-interop runtime)
# always materializes a distinct <w:r> for a sub-range once its run
# properties are touched, even if the net formatting ends up identical to
# the surrounding text.
function Split-Run($lo, $hi, $needle) {
    $rng = $d.Range($lo, $hi)
    $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $rng.Bold = 1
    $rng.Bold = 0
}

# --- Paragraph: "Наличие описанного извне шаблона структуры name с двумя ..." ---
$p = $d.Paragraphs.Item(13)
Split-Run $p.Range.Start $p.Range.End "name"

# --- Paragraph: "Наличие описанного извне шаблона структуры  «пассажир» с элементами: структурой name, ном_рейса (строка), кол_вещей(целый тип), общ_вес(целый тип)." ---
$p = $d.Paragraphs.Item(14)
Split-Run $p.Range.Start $p.Range.End "name"
Split-Run $p.Range.Start $p.Range.End "ном_рейса"
Split-Run $p.Range.Start $p.Range.End "кол_вещей"
Split-Run $p.Range.Start $p.Range.End "общ_вес"

# --- Paragraph: "Функция main() должна объявлять массив структур ..." ---
$p = $d.Paragraphs.Item(15)
Split-Run $p.Range.Start $p.Range.End "main"

# --- "Анализ требований" heading gains " и разработка алгоритма" ---
$p = $d.Paragraphs.Item(20)
$p.Range.InsertAfter(" и разработка алгоритма")
# force that appended text to live in its own run, matching the diff's
# two-run paragraph
Split-Run $p.Range.Start $p.Range.End " и разработка алгоритма"

# --- Remove the old "Перед тем ..." paragraph through the "Функции" list of
#     declarations, keeping only a single empty paragraph (which carries the
#     _GoBack bookmark) in their place. ---
$first = $d.Paragraphs.Item(21)
$last = $d.Paragraphs.Item(29)
$rng = $d.Range($first.Range.Start, $last.Range.End)
$rng.Text = ""

# Re-anchor the _GoBack bookmark onto the now-empty paragraph left behind.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $d.Paragraphs.Item(21).Range)
